# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# Both sheets carry the same event rows, so the same updates apply to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 118
    8  = 452
    13 = 297
    33 = 271
    35 = 51
    40 = 3610
    43 = 906
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
